{"js": "// Edit: trim the trailing clause \", until Smartcash reaches a considerable\n// market cap\" off the ASICs/Keccak-mining sentence, so it now ends\n// \"...quite some time.\"\n//\n// Before: \"...no ASICs will be created for quite some time, until Smartcash\n//          reaches a considerable market cap.\"\n// After:  \"...no ASICs will be created for quite some time.\"\n\nconst body = context.document.body;\n\nconst clauseToRemove = \", until Smartcash reaches a considerable market cap\";\n\nconst searchResults = body.search(clauseToRemove, { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  // Delete just the matched clause, leaving the final period untouched so\n  // the surrounding run's formatting (rPr) is preserved.\n  searchResults.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Edit: trim the trailing clause \", until Smartcash reaches a considerable\n# market cap\" off the ASICs/Keccak-mining sentence, so it now ends\n# \"...quite some time.\"\n#\n# Before: \"...no ASICs will be created for quite some time, until Smartcash\n#          reaches a considerable market cap.\"\n# After:  \"...no ASICs will be created for quite some time.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n    \", until Smartcash reaches a considerable market cap\",  # FindText (leave the trailing period in place)\n    $true,    # MatchCase\n    $false,   # MatchWholeWord\n    $false,   # MatchWildcards\n    $false,   # MatchSoundsLike\n    $false,   # MatchAllWordForms\n    $true,    # Forward\n    1,        # Wrap (wdFindContinue)\n    $false,   # Format\n    \"\",       # ReplaceWith (delete the clause)\n    2         # Replace (wdReplaceAll)\n) | Out-Null\n"}
